$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at sheet row 65 (pushes the existing rows 65:90 down to 66:91)
# so the new "Pedro Neira" / Bio Bio candidate lands right after the other
# Bio Bio rows and before the Araucania group.
$ws.Rows("65").Insert()

# Fill in the new row's values (column A is the running id = row-1, same
# pattern as every other row in the sheet).
$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = 9
$ws.Cells.Item(65, 3).Value = 8
$ws.Cells.Item(65, 4).Value = "Bío Bío"
$ws.Cells.Item(65, 5).Value = "Pedro Neira"
$ws.Cells.Item(65, 6).Value = "Humanicemos Chile"
$ws.Cells.Item(65, 7).Value = "PH"
$ws.Cells.Item(65, 8).Value = 0

# Update the selection / scroll position to match the saved view.
$ws.Range("A3:A91").Select()
$excel.ActiveWindow.ScrollRow = 72
